$wb = $excel.ActiveWorkbook

# --- New test-data values in the "filter value" column (C3) of each sheet ---
# ConclToApprove gets its own value; the other four sheets share the other value
# (matches the order new shared strings were introduced: "124$" then "123$").
$wsConcl = $wb.Worksheets.Item("ConclToApprove")
$wsConcl.Range("C3").Value = "124$"

$wsDonation = $wb.Worksheets.Item("DonationInfo")
$wsDonation.Range("C3").Value = "123$"

$wsTest = $wb.Worksheets.Item("TestInfo")
$wsTest.Range("C3").Value = "123$"

$wsWorklistDetail = $wb.Worksheets.Item("WorklistDetail")
$wsWorklistDetail.Range("C3").Value = "123$"

$wsWorklistOverview = $wb.Worksheets.Item("WorklistOverview")
$wsWorklistOverview.Range("C3").Value = "123$"

# --- The active/selected tab moves from WorklistOverview to ConclToApprove, ---
# --- with the new selection sitting on E8 (WorklistOverview keeps its G17 selection). ---
$wsConcl.Activate()
$wsConcl.Range("E8").Select()
